$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at the top of the "Agrícola del Norte S.A. de Arica - Piña"
# weekly price block (rows 222-225), shifting the existing rows 222-284 down
# to 226-288. Excel copies the row-above formatting automatically (preserves
# the date style on column D, etc.)
$ws.Rows("222:225").Insert()

# Populate the 4 newly inserted rows with the new week's data (same market /
# product metadata as before, updated report date and updated volume/price
# figures for the "Primera", "Segunda" and "Tercera" quality grades).

# Row 222 - Especial
$ws.Cells.Item(222, 1).Value = 1
$ws.Cells.Item(222, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(222, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(222, 4).Value = 44985
$ws.Cells.Item(222, 5).Value = 15
$ws.Cells.Item(222, 6).Value = "Fruta"
$ws.Cells.Item(222, 7).Value = 100108
$ws.Cells.Item(222, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(222, 9).Value = 100108005
$ws.Cells.Item(222, 10).Value = "Piña"
$ws.Cells.Item(222, 11).Value = "Caramelo"
$ws.Cells.Item(222, 12).Value = "Especial"
$ws.Cells.Item(222, 13).Value = 200
$ws.Cells.Item(222, 14).Value = 20000
$ws.Cells.Item(222, 15).Value = 21000
$ws.Cells.Item(222, 16).Value = 20500
$ws.Cells.Item(222, 17).Value = "$/caja 10 unidades"
$ws.Cells.Item(222, 18).Value = "Ecuador"
$ws.Cells.Item(222, 19).Value = 2050
$ws.Cells.Item(222, 20).Value = 10

# Row 223 - Primera
$ws.Cells.Item(223, 1).Value = 1
$ws.Cells.Item(223, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(223, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(223, 4).Value = 44985
$ws.Cells.Item(223, 5).Value = 15
$ws.Cells.Item(223, 6).Value = "Fruta"
$ws.Cells.Item(223, 7).Value = 100108
$ws.Cells.Item(223, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(223, 9).Value = 100108005
$ws.Cells.Item(223, 10).Value = "Piña"
$ws.Cells.Item(223, 11).Value = "Caramelo"
$ws.Cells.Item(223, 12).Value = "Primera"
$ws.Cells.Item(223, 13).Value = 190
$ws.Cells.Item(223, 14).Value = 20000
$ws.Cells.Item(223, 15).Value = 21000
$ws.Cells.Item(223, 16).Value = 20474
$ws.Cells.Item(223, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(223, 18).Value = "Ecuador"
$ws.Cells.Item(223, 19).Value = 1706
$ws.Cells.Item(223, 20).Value = 12

# Row 224 - Segunda
$ws.Cells.Item(224, 1).Value = 1
$ws.Cells.Item(224, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(224, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(224, 4).Value = 44985
$ws.Cells.Item(224, 5).Value = 15
$ws.Cells.Item(224, 6).Value = "Fruta"
$ws.Cells.Item(224, 7).Value = 100108
$ws.Cells.Item(224, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(224, 9).Value = 100108005
$ws.Cells.Item(224, 10).Value = "Piña"
$ws.Cells.Item(224, 11).Value = "Caramelo"
$ws.Cells.Item(224, 12).Value = "Segunda"
$ws.Cells.Item(224, 13).Value = 220
$ws.Cells.Item(224, 14).Value = 20000
$ws.Cells.Item(224, 15).Value = 21000
$ws.Cells.Item(224, 16).Value = 20545
$ws.Cells.Item(224, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(224, 18).Value = "Ecuador"
$ws.Cells.Item(224, 19).Value = 1468
$ws.Cells.Item(224, 20).Value = 14

# Row 225 - Tercera
$ws.Cells.Item(225, 1).Value = 1
$ws.Cells.Item(225, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(225, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(225, 4).Value = 44985
$ws.Cells.Item(225, 5).Value = 15
$ws.Cells.Item(225, 6).Value = "Fruta"
$ws.Cells.Item(225, 7).Value = 100108
$ws.Cells.Item(225, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(225, 9).Value = 100108005
$ws.Cells.Item(225, 10).Value = "Piña"
$ws.Cells.Item(225, 11).Value = "Caramelo"
$ws.Cells.Item(225, 12).Value = "Tercera"
$ws.Cells.Item(225, 13).Value = 250
$ws.Cells.Item(225, 14).Value = 20000
$ws.Cells.Item(225, 15).Value = 21000
$ws.Cells.Item(225, 16).Value = 20400
$ws.Cells.Item(225, 17).Value = "$/caja 16 unidades"
$ws.Cells.Item(225, 18).Value = "Ecuador"
$ws.Cells.Item(225, 19).Value = 1275
$ws.Cells.Item(225, 20).Value = 16
